# edit.ps1 - apply SMARTNODES.docx (Italian) edits via Word COM-interop
#
# Summary of the content changes performed here (see task diff):
#   1. The bullet "InstantPay (bold) (instant transactions). In contrast,
#      Bitcoin takes about 10 minutes to confirm a payment." is replaced
#      with a single, unbolded run of new copy about InstantPay, using
#      the "Open Sans"/#252525 styling and the numId=2 (tighter) list
#      formatting/spacing that the rest of the post-edit document uses.
#   2. The bullet "Will have more services added later" is replaced with
#      a new SmartRewards description, again switching to the
#      "Open Sans"/#252525 styling and numId=2 spacing.
#
$d = $word.ActiveDocument

function Find-ParagraphIndexByText($needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

function Set-ParagraphXml($paraIndex, $pPrXml, $runsXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $full = $p.Range
    $xmlSnippet = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $pPrXml + $runsXml + '</w:p><w:p></w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $full.InsertXML($xmlSnippet)
}

# -- Paragraph 1: "InstantPay ... ten minutes to confirm a payment." ------
$instantPayIdx = Find-ParagraphIndexByText "instant transactions"
if ($instantPayIdx -eq -1) {
    $instantPayIdx = Find-ParagraphIndexByText "InstantPay"
}

$pPr1 = '<w:pPr>' +
            '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' +
            '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
            '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
            '<w:textAlignment w:val="baseline"/>' +
            '<w:rPr>' +
                '<w:rFonts w:ascii="Open Sans" w:eastAsia="Times New Roman" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
                '<w:color w:val="252525"/>' +
                '<w:sz w:val="21"/>' +
                '<w:szCs w:val="21"/>' +
            '</w:rPr>' +
        '</w:pPr>'
$run1 = '<w:r>' +
            '<w:rPr>' +
                '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/>' +
                '<w:color w:val="252525"/>' +
                '<w:sz w:val="21"/>' +
                '<w:szCs w:val="21"/>' +
            '</w:rPr>' +
            '<w:t xml:space="preserve">InstantPay (Instant Transactions): Allows for SmartCash transactions to be locked in about a second. No risk of double spending a transaction, so the receiver can trust that transaction immediately.</w:t>' +
        '</w:r>'

Set-ParagraphXml $instantPayIdx $pPr1 $run1

# -- Paragraph 2: "Will have more services added later" -------------------
$smartRewardsIdx = Find-ParagraphIndexByText "Will have more services"

$pPr2 = '<w:pPr>' +
            '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' +
            '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
            '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
            '<w:textAlignment w:val="baseline"/>' +
            '<w:rPr>' +
                '<w:rFonts w:ascii="Open Sans" w:eastAsia="Times New Roman" w:hAnsi="Open Sans" w:cs="Open Sans"/>' +
                '<w:color w:val="252525"/>' +
                '<w:sz w:val="21"/>' +
                '<w:szCs w:val="21"/>' +
            '</w:rPr>' +
        '</w:pPr>'
$run2 = '<w:r>' +
            '<w:rPr>' +
                '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/>' +
                '<w:color w:val="252525"/>' +
                '<w:sz w:val="21"/>' +
                '<w:szCs w:val="21"/>' +
            '</w:rPr>' +
            '<w:t xml:space="preserve">SmartRewards: SmartRewards are calculated by the SmartNodes to allow for distribution to be handled automatically by the block rewards.</w:t>' +
        '</w:r>'

Set-ParagraphXml $smartRewardsIdx $pPr2 $run2

Write-Output "Edits applied: InstantPay paragraph index=$instantPayIdx, SmartRewards paragraph index=$smartRewardsIdx"
